$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "57.578.65"
$ws.Cells.Item(2, 5).Value = "  -4.74%  "
$ws.Cells.Item(3, 4).Value = "3.148.18"
$ws.Cells.Item(3, 5).Value = "  -5.30%  "
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 4).Value = "'526.40"
$ws.Cells.Item(5, 5).Value = "  -6.56%  "
$ws.Cells.Item(6, 4).Value = "'133.50"
$ws.Cells.Item(6, 5).Value = "  -8.17%  "
$ws.Cells.Item(7, 5).Value = "  -0.14%  "
$ws.Cells.Item(8, 4).Value = "3.145.74"
$ws.Cells.Item(8, 5).Value = "  -5.32%  "
$ws.Cells.Item(9, 4).Value = "'0.449"
$ws.Cells.Item(9, 5).Value = "  -6.98%  "
$ws.Cells.Item(10, 4).Value = "'7.23"
$ws.Cells.Item(10, 5).Value = "  -8.01%  "
$ws.Cells.Item(11, 4).Value = "'0.111"
$ws.Cells.Item(11, 5).Value = "  -8.05%  "
$ws.Cells.Item(12, 4).Value = "'0.392"
$ws.Cells.Item(12, 5).Value = "  -4.75%  "
$ws.Cells.Item(13, 4).Value = "3.689.74"
$ws.Cells.Item(13, 5).Value = "  -5.42%  "
$ws.Cells.Item(14, 4).Value = "'0.127"
$ws.Cells.Item(14, 5).Value = "  -1.85%  "
$ws.Cells.Item(15, 4).Value = "'25.85"
$ws.Cells.Item(15, 5).Value = "  -6.16%  "
$ws.Cells.Item(16, 4).Value = "3.150.13"
$ws.Cells.Item(16, 5).Value = "  -5.90%  "
$ws.Cells.Item(17, 4).Value = "57.520.64"
$ws.Cells.Item(17, 5).Value = "  -4.98%  "
$ws.Cells.Item(18, 4).Value = "'0.0000153"
$ws.Cells.Item(18, 5).Value = "  -8.65%  "
$ws.Cells.Item(19, 4).Value = "'5.82"
$ws.Cells.Item(19, 5).Value = "  -6.52%  "
$ws.Cells.Item(20, 4).Value = "'13.09"
$ws.Cells.Item(20, 5).Value = "  -9.31%  "
$ws.Cells.Item(21, 4).Value = "'8.02"
$ws.Cells.Item(21, 5).Value = "  -9.16%  "
$ws.Cells.Item(22, 4).Value = "'346.52"
$ws.Cells.Item(22, 5).Value = "  -7.72%  "
$ws.Cells.Item(23, 4).Value = "'0.999"
$ws.Cells.Item(23, 5).Value = "  -0.13%  "
$ws.Cells.Item(24, 4).Value = "'69.58"
$ws.Cells.Item(24, 5).Value = "  -6.83%  "
$ws.Cells.Item(25, 4).Value = "'0.511"
$ws.Cells.Item(25, 5).Value = "  -7.88%  "
$ws.Cells.Item(26, 4).Value = "3.290.08"
$ws.Cells.Item(26, 5).Value = "  -5.63%  "
$ws.Cells.Item(27, 4).Value = "0.0₃0962"
$ws.Cells.Item(27, 5).Value = "  -9.88%  "
$ws.Cells.Item(28, 5).Value = "  -4.03%  "
$ws.Cells.Item(29, 4).Value = "'1.00"
$ws.Cells.Item(29, 5).Value = "  -0.12%  "
$ws.Cells.Item(30, 4).Value = "'6.86"
$ws.Cells.Item(30, 5).Value = "  -5.91%  "
$ws.Cells.Item(31, 4).Value = "'0.996"
$ws.Cells.Item(31, 5).Value = "  -0.49%  "
$ws.Cells.Item(32, 4).Value = "'1.88"
$ws.Cells.Item(32, 5).Value = "  -9.35%  "
$ws.Cells.Item(33, 4).Value = "'6.94"
$ws.Cells.Item(33, 5).Value = "  -9.19%  "
$ws.Cells.Item(34, 4).Value = "'21.60"
$ws.Cells.Item(34, 5).Value = "  -5.01%  "
$ws.Cells.Item(35, 5).Value = "  -5.79%  "
$ws.Cells.Item(36, 4).Value = "'4.97"
$ws.Cells.Item(36, 5).Value = "  -5.20%  "
$ws.Cells.Item(37, 4).Value = "'158.76"
$ws.Cells.Item(37, 5).Value = "  -4.80%  "
$ws.Cells.Item(38, 4).Value = "'6.24"
$ws.Cells.Item(38, 5).Value = "  -8.05%  "
$ws.Cells.Item(39, 4).Value = "'1.41"
$ws.Cells.Item(39, 5).Value = "  -8.65%  "
$ws.Cells.Item(40, 4).Value = "'25.95"
$ws.Cells.Item(40, 5).Value = "  -6.96%  "
$ws.Cells.Item(41, 4).Value = "'0.0697"
$ws.Cells.Item(41, 5).Value = "  -5.87%  "
$ws.Cells.Item(42, 4).Value = "3.174.73"
$ws.Cells.Item(42, 5).Value = "  -5.78%  "
$ws.Cells.Item(43, 4).Value = "'40.27"
$ws.Cells.Item(43, 5).Value = "  -4.36%  "
$ws.Cells.Item(44, 4).Value = "'0.694"
$ws.Cells.Item(44, 5).Value = "  -7.86%  "
$ws.Cells.Item(45, 4).Value = "'1.08"
$ws.Cells.Item(45, 5).Value = "  -4.74%  "
$ws.Cells.Item(46, 4).Value = "'3.96"
$ws.Cells.Item(47, 4).Value = "'0.999"
$ws.Cells.Item(47, 5).Value = "  -0.26%  "
$ws.Cells.Item(48, 4).Value = "'1.46"
$ws.Cells.Item(48, 5).Value = "  -8.54%  "
$ws.Cells.Item(49, 4).Value = "2.268.36"
$ws.Cells.Item(49, 5).Value = "  -7.14%  "
$ws.Cells.Item(50, 4).Value = "'6.24"
$ws.Cells.Item(50, 5).Value = "  -5.95%  "
$ws.Cells.Item(51, 4).Value = "'20.61"
$ws.Cells.Item(51, 5).Value = "  -6.84%  "
